$d = $word.ActiveDocument

# 1. Merge "Erstellung der " + "Prak" + " Website" (with spell-check
#    proofErr markup wrapped around "Prak") into a single plain run.
#    Find/Replace across the whole story collapses the matched runs
#    (and drops the proofErr markers) into one run holding the replacement
#    text, which is exactly what the diff shows.
$d.Content.Find.Execute(
    "Erstellung der Prak Website", $true, $false, $false, $false, $false,
    $true, 1, $false, "Erstellung der Prak Website", 2) | Out-Null

# 2. Merge the "Fake" run (and its proofErr wrapper) back into the
#    surrounding sentence as a single run, same technique as above.
$old2 = " Es wurde sich kollektiv entschieden das einzelne Schüler nicht in der Lage sein sollten Unternehmen selbst einzutragen. Der Grund dafür ist das der Aufwand höher ist zu kontrollieren ob irgendein Fake Unternehmen eingetragen wurde, als einem Lehrer die Daten über das Unternehmen zu geben."
$d.Content.Find.Execute(
    $old2, $true, $false, $false, $false, $false,
    $true, 1, $false, $old2, 2) | Out-Null

# 3. Replace the lone "Sebastian" paragraph (under the "Probleme" heading)
#    with the new problem description text. Scope the Find to that single
#    paragraph's range so the other "Sebastian" occurrences elsewhere in
#    the document are left untouched. Remember its index so step 4 below
#    only touches the blank paragraphs that immediately follow it.
$sebastianIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Sebastian`r") {
        $r = $p.Range
        $r.End = $r.End - 1
        $r.Find.Execute(
            "Sebastian", $true, $false, $false, $false, $false,
            $true, 1, $false, "Das Hauptproblem war der Import der  Datenbank ", 2) | Out-Null
        $sebastianIndex = $i
        break
    }
}

# 4. Remove one of the two consecutive empty paragraphs that followed that
#    paragraph, collapsing the double blank line to a single blank line.
if ($sebastianIndex -ge 1) {
    for ($i = $sebastianIndex; $i -le $d.Paragraphs.Count - 1; $i++) {
        $p1 = $d.Paragraphs.Item($i)
        $p2 = $d.Paragraphs.Item($i + 1)
        if ($p1.Range.Text -eq "`r" -and $p2.Range.Text -eq "`r") {
            $p1.Range.Delete() | Out-Null
            break
        }
    }
}
